$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Cell, $Val)
    $rng = $ws.Range($Cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $Val
    $rng.Style = $origStyle
}

Set-TextValue "D2" '37.737.78'
Set-TextValue "E2" '  -0.17%  '
Set-TextValue "D3" '2.038.12'
Set-TextValue "E3" '  +0.49%  '
Set-TextValue "E4" '  -0.02%  '
Set-TextValue "D5" '227.15'
Set-TextValue "E5" '  +0.06%  '
Set-TextValue "D6" '0.606'
Set-TextValue "E6" '  -0.87%  '
Set-TextValue "D7" '59.45'
Set-TextValue "E7" '  -0.29%  '
Set-TextValue "E8" '  +0.03%  '
Set-TextValue "E9" '  -2.39%  '
Set-TextValue "D10" '0.0834'
Set-TextValue "E10" '  +3.06%  '
Set-TextValue "E11" '  -0.16%  '
Set-TextValue "D12" '2.340.37'
Set-TextValue "E12" '  +0.53%  '
Set-TextValue "D13" '14.43'
Set-TextValue "E13" '  -0.77%  '
Set-TextValue "D14" '21.04'
Set-TextValue "E14" '  +0.27%  '
Set-TextValue "D15" '5.44'
Set-TextValue "E15" '  +4.24%  '
Set-TextValue "E16" '  +2.59%  '
Set-TextValue "D17" '2.025.30'
Set-TextValue "E17" '  -1.30%  '
Set-TextValue "D18" '37.731.33'
Set-TextValue "E18" '  -0.24%  '
Set-TextValue "D19" '5.93'
Set-TextValue "E19" '  -1.56%  '
Set-TextValue "D20" '69.36'
Set-TextValue "E20" '  -0.18%  '
Set-TextValue "D21" '0.0₃0822'
Set-TextValue "E21" '  +0.03%  '
Set-TextValue "D22" '223.83'
Set-TextValue "E22" '  -0.35%  '
Set-TextValue "E23" '  +0.05%  '
Set-TextValue "D24" '2.43'
Set-TextValue "E24" '  +0.70%  '
Set-TextValue "D25" '2.27'
Set-TextValue "E25" '  +2.44%  '
Set-TextValue "D26" '167.93'
Set-TextValue "E26" '  +1.74%  '
Set-TextValue "D27" '9.36'
Set-TextValue "E27" '  +2.16%  '
Set-TextValue "E28" '  -1.03%  '
Set-TextValue "D29" '18.77'
Set-TextValue "E29" '  -0.45%  '
Set-TextValue "E30" '  -0.38%  '
Set-TextValue "E31" '  +0.07%  '
Set-TextValue "D32" '2.23'
Set-TextValue "E32" '  +9.46%  '
Set-TextValue "D34" '0.0607'
Set-TextValue "E34" '  +1.10%  '
Set-TextValue "D35" '4.47'
Set-TextValue "E35" '  -0.32%  '
Set-TextValue "D36" '6.52'
Set-TextValue "E36" '  +3.05%  '
Set-TextValue "E37" '  +4.23%  '
Set-TextValue "E38" '  +5.26%  '
Set-TextValue "D39" '1.00'
Set-TextValue "E39" '  -0.18%  '
Set-TextValue "D40" '18.06'
Set-TextValue "E40" '  +9.03%  '
Set-TextValue "D41" '1.530.68'
Set-TextValue "E41" '  -0.66%  '
Set-TextValue "D42" '97.07'
Set-TextValue "E42" '  +0.72%  '
Set-TextValue "E43" '  -0.58%  '
Set-TextValue "B44" 'FTXToken'
Set-TextValue "C44" 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue "D44" '4.36'
Set-TextValue "E44" '  +11.55%  '
Set-TextValue "B45" 'HuobiToken'
Set-TextValue "C45" 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue "D45" '2.84'
Set-TextValue "E45" '  +0.61%  '
Set-TextValue "D46" '0.0906'
Set-TextValue "E46" '  -1.42%  '
Set-TextValue "E47" '  +0.64%  '
Set-TextValue "E48" '  +0.48%  '
Set-TextValue "E49" '  -0.96%  '
Set-TextValue "E50" '  -0.31%  '
Set-TextValue "D51" '2.230.19'
Set-TextValue "E51" '  +0.54%  '
